# Insert two new data rows at row 12 (pushing existing rows 12-99 down to 14-101)
# and populate them with the new weekly price records ("Fruta / hortaliza, semanal").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows starting at row 12.
$ws.Rows.Item(12).Resize(2).Insert()

# Constant values shared by every data row in this sheet.
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$codreg    = 15
$catId     = 100112028
$categoria = "Sandia"
$variedad  = "Sin especificar"
$unidad    = "$/kilo (volumen en unidades)"
$clasif    = "Hortaliza"

# New row 12 data.
$ws.Cells.Item(12, 1).Value  = 1
$ws.Cells.Item(12, 2).Value  = $mercado
$ws.Cells.Item(12, 3).Value  = $region
$ws.Cells.Item(12, 4).Value  = 45282
$ws.Cells.Item(12, 5).Value  = $codreg
$ws.Cells.Item(12, 6).Value  = $catId
$ws.Cells.Item(12, 7).Value  = $categoria
$ws.Cells.Item(12, 8).Value  = $variedad
$ws.Cells.Item(12, 9).Value  = "Extra"
$ws.Cells.Item(12, 10).Value = 500
$ws.Cells.Item(12, 11).Value = 450
$ws.Cells.Item(12, 12).Value = 500
$ws.Cells.Item(12, 13).Value = 475
$ws.Cells.Item(12, 14).Value = $unidad
$ws.Cells.Item(12, 15).Value = "Perú"
$ws.Cells.Item(12, 16).Value = 475
$ws.Cells.Item(12, 17).Value = 1
$ws.Cells.Item(12, 18).Value = $clasif

# New row 13 data.
$ws.Cells.Item(13, 1).Value  = 1
$ws.Cells.Item(13, 2).Value  = $mercado
$ws.Cells.Item(13, 3).Value  = $region
$ws.Cells.Item(13, 4).Value  = 45282
$ws.Cells.Item(13, 5).Value  = $codreg
$ws.Cells.Item(13, 6).Value  = $catId
$ws.Cells.Item(13, 7).Value  = $categoria
$ws.Cells.Item(13, 8).Value  = $variedad
$ws.Cells.Item(13, 9).Value  = "Primera"
$ws.Cells.Item(13, 10).Value = 600
$ws.Cells.Item(13, 11).Value = 450
$ws.Cells.Item(13, 12).Value = 500
$ws.Cells.Item(13, 13).Value = 471
$ws.Cells.Item(13, 14).Value = $unidad
$ws.Cells.Item(13, 15).Value = "Perú"
$ws.Cells.Item(13, 16).Value = 471
$ws.Cells.Item(13, 17).Value = 1
$ws.Cells.Item(13, 18).Value = $clasif

# Ensure the date cells use the same date/time number format as the rest of column D.
$ws.Cells.Item(12, 4).NumberFormat = $ws.Cells.Item(14, 4).NumberFormat
$ws.Cells.Item(13, 4).NumberFormat = $ws.Cells.Item(14, 4).NumberFormat
